# EnterpriseJavaTimeLog.xlsx edit script
# Adds separate error displays for 403's and other errors; logs Week 7 progress.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 4 blank rows before row 52 (pushes existing rows 52+ down to 56+) ---
$ws.Rows("52:55").Insert()

# --- 2. Row 43: hours corrected from 0.5 to 1 ---
$ws.Range("B43").Value = 1

# --- 3. Row 44: new Week 7 log entry (date, hours, task) ---
$ws.Range("A44").Value = 43540
$ws.Range("B44").Value = 2.5
$ws.Range("D44").Value = "Week 7: Completed readings (light treatment) and videos`nIndie Project: added separate jsp's for 403 errors and other errors"
$ws.Rows(44).RowHeight = 30

# --- 4. Row 47: remove the "Thurs..." note (it moves down to row 51) ---
$ws.Range("D47").Clear()

# --- 5. Row 48: new TODO note about custom error pages ---
$ws.Range("D48").Value = "TODO WEEK 7 - • TODO google “web xml custom error pages” and implement`n• Maybe this one https://www.tutorialspoint.com/servlets/servlets-exception-handling.htm`n"
$ws.Rows(48).RowHeight = 45

# --- 6. Row 51: the "Thurs..." note now lives here ---
$ws.Range("D51").Value = "Thurs -  a little more time than listed."

# --- 7. Fix up the sheet view: scroll position and active selection ---
$ws.Range("D45").Select()
